# Fixes for bugs #98 and #101
#
# Bug #98: The "RunInParallel" setting on the Config sheet was left as
#          "No" - flip it to "Yes".
# Bug #101: Test case 101's Execute column on the "Test Cases" sheet held
#           a stray "TestCaseNumber=101" value instead of the intended
#           "Groups=Regression" filter expression - correct it.

$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("Config")
$wsTests  = $wb.Worksheets.Item("Test Cases")

# Bug #98 - RunInParallel should be enabled.
$wsConfig.Range("B3").Value = "Yes"

# Bug #101 - fix the Execute expression for test case 101.
$wsTests.Range("D2").Value = "Groups=Regression"

# Leave the UI focused back on the Config sheet (selection on B10), with
# the Test Cases sheet's own selection parked on D3.
$wsTests.Range("D3").Select()
$wsConfig.Activate()
$wsConfig.Range("B10").Select()
